$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new "hierarchy" / "Hierarchy" key/value row just before the
#     existing "instruct_collect" row (old row 76) ---
$ws.Rows("76:76").Insert()
$ws.Range("A76").Value = "hierarchy"
$ws.Range("B76").Value = "Hierarchy"

# --- Insert new "instruct_category_detail" key/value row right after
#     "instruct_category_drag" (old row 77, now row 78) ---
$ws.Rows("79:79").Insert()
$ws.Range("A79").Value = "instruct_category_detail"
$ws.Range("B79").Value = "Press any of the categories to review them."

# --- Update the level title text values to account for the newly
#     inserted "3 - Polygons #3" level (everything renumbered by one) ---
$ws.Range("B107").Value = "3 - Polygons #3"
$ws.Range("B108").Value = "4 - Triangles (Angles)"
$ws.Range("B109").Value = "5 - Triangles (Sides)"
$ws.Range("B110").Value = "6 - Triangles (All)"
$ws.Range("B111").Value = "7 - Quads (Parallelograms)"
$ws.Range("B112").Value = "8 - Quads (Types)"

# --- Append a new row for the "level_title_8" key (the hierarchy review
#     level title) at the end of the table ---
$ws.Rows("113:113").Insert()
$ws.Range("A113").Value = "level_title_8"
$ws.Range("B113").Value = "9 - Quads (Hierarchy)"

# --- Update the active selection / view to rest on the newly inserted
#     row, matching the saved workbook view state ---
$ws.Range("A76").Select()
